# Update column F ("想去人数" / number of people interested) values on the
# "展览" and "全部类型" sheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F value
$updates = @{
    2  = 1579
    4  = 1037
    5  = 33
    6  = 69
    7  = 2722
    9  = 1745
    12 = 596
    14 = 18
    15 = 129
    16 = 80
    18 = 18
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
